$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (2-7) first, then rewrite rows 2-10 with updated values.
$ws.Range("A2:T7").ClearContents()

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Alcam"
$ws.Cells.Item(2, 3).Value = "L1cam"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 42.078429
$ws.Cells.Item(2, 8).Value = 126.235287
$ws.Cells.Item(2, 9).Value = 0.9777085588730982
$ws.Cells.Item(2, 10).Value = 0.9777085588730982
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 22.59487733333333
$ws.Cells.Item(2, 14).Value = 67.784632
$ws.Cells.Item(2, 15).Value = 0.7395019553569895
$ws.Cells.Item(2, 16).Value = 0.7395019553569895
$ws.Cells.Item(2, 17).Value = 950.756941634376
$ws.Cells.Item(2, 18).Value = 8556.812474709384
$ws.Cells.Item(2, 19).Value = 0.7230173910559203
$ws.Cells.Item(2, 20).Value = 0.7230173910559203

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Alcam"
$ws.Cells.Item(3, 3).Value = "L1cam"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 42.078429
$ws.Cells.Item(3, 8).Value = 126.235287
$ws.Cells.Item(3, 9).Value = 0.9777085588730982
$ws.Cells.Item(3, 10).Value = 0.9777085588730982
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3045986666666667
$ws.Cells.Item(3, 14).Value = 0.913796
$ws.Cells.Item(3, 15).Value = 0.00996913177602551
$ws.Cells.Item(3, 16).Value = 0.00996913177602551
$ws.Cells.Item(3, 17).Value = 12.817033368828
$ws.Cells.Item(3, 18).Value = 115.353300319452
$ws.Cells.Item(3, 19).Value = 0.009746905461953911
$ws.Cells.Item(3, 20).Value = 0.009746905461953911

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Alcam"
$ws.Cells.Item(4, 3).Value = "L1cam"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 42.078429
$ws.Cells.Item(4, 8).Value = 126.235287
$ws.Cells.Item(4, 9).Value = 0.9777085588730982
$ws.Cells.Item(4, 10).Value = 0.9777085588730982
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.654706
$ws.Cells.Item(4, 14).Value = 22.964118
$ws.Cells.Item(4, 15).Value = 0.2505289128669849
$ws.Cells.Item(4, 16).Value = 0.2505289128669849
$ws.Cells.Item(4, 17).Value = 322.098002936874
$ws.Cells.Item(4, 18).Value = 2898.882026431866
$ws.Cells.Item(4, 19).Value = 0.2449442623552238
$ws.Cells.Item(4, 20).Value = 0.2449442623552238

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Alcam"
$ws.Cells.Item(5, 3).Value = "L1cam"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.4771596666666666
$ws.Cells.Item(5, 8).Value = 1.431479
$ws.Cells.Item(5, 9).Value = 0.01108698925164327
$ws.Cells.Item(5, 10).Value = 0.01108698925164327
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 22.59487733333333
$ws.Cells.Item(5, 14).Value = 67.784632
$ws.Cells.Item(5, 15).Value = 0.7395019553569895
$ws.Cells.Item(5, 16).Value = 0.7395019553569895
$ws.Cells.Item(5, 17).Value = 10.78136413674756
$ws.Cells.Item(5, 18).Value = 97.032277230728
$ws.Cells.Item(5, 19).Value = 0.008198850230612124
$ws.Cells.Item(5, 20).Value = 0.008198850230612124

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Alcam"
$ws.Cells.Item(6, 3).Value = "L1cam"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.4771596666666666
$ws.Cells.Item(6, 8).Value = 1.431479
$ws.Cells.Item(6, 9).Value = 0.01108698925164327
$ws.Cells.Item(6, 10).Value = 0.01108698925164327
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.3045986666666667
$ws.Cells.Item(6, 14).Value = 0.913796
$ws.Cells.Item(6, 15).Value = 0.00996913177602551
$ws.Cells.Item(6, 16).Value = 0.00996913177602551
$ws.Cells.Item(6, 17).Value = 0.1453421982537778
$ws.Cells.Item(6, 18).Value = 1.308079784284
$ws.Cells.Item(6, 19).Value = 0.0001105276568490102
$ws.Cells.Item(6, 20).Value = 0.0001105276568490102

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Alcam"
$ws.Cells.Item(7, 3).Value = "L1cam"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.4771596666666666
$ws.Cells.Item(7, 8).Value = 1.431479
$ws.Cells.Item(7, 9).Value = 0.01108698925164327
$ws.Cells.Item(7, 10).Value = 0.01108698925164327
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.654706
$ws.Cells.Item(7, 14).Value = 22.964118
$ws.Cells.Item(7, 15).Value = 0.2505289128669849
$ws.Cells.Item(7, 16).Value = 0.2505289128669849
$ws.Cells.Item(7, 17).Value = 3.652516963391333
$ws.Cells.Item(7, 18).Value = 32.872652670522
$ws.Cells.Item(7, 19).Value = 0.002777611364182135
$ws.Cells.Item(7, 20).Value = 0.002777611364182135

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Alcam"
$ws.Cells.Item(8, 3).Value = "L1cam"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.482215
$ws.Cells.Item(8, 8).Value = 1.446645
$ws.Cells.Item(8, 9).Value = 0.01120445187525872
$ws.Cells.Item(8, 10).Value = 0.01120445187525872
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 22.59487733333333
$ws.Cells.Item(8, 14).Value = 67.784632
$ws.Cells.Item(8, 15).Value = 0.7395019553569895
$ws.Cells.Item(8, 16).Value = 0.7395019553569895
$ws.Cells.Item(8, 17).Value = 10.89558877329333
$ws.Cells.Item(8, 18).Value = 98.06029895964
$ws.Cells.Item(8, 19).Value = 0.008285714070457112
$ws.Cells.Item(8, 20).Value = 0.008285714070457112

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Alcam"
$ws.Cells.Item(9, 3).Value = "L1cam"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.482215
$ws.Cells.Item(9, 8).Value = 1.446645
$ws.Cells.Item(9, 9).Value = 0.01120445187525872
$ws.Cells.Item(9, 10).Value = 0.01120445187525872
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.3045986666666667
$ws.Cells.Item(9, 14).Value = 0.913796
$ws.Cells.Item(9, 15).Value = 0.00996913177602551
$ws.Cells.Item(9, 16).Value = 0.00996913177602551
$ws.Cells.Item(9, 17).Value = 0.1468820460466667
$ws.Cells.Item(9, 18).Value = 1.32193841442
$ws.Cells.Item(9, 19).Value = 0.0001116986572225903
$ws.Cells.Item(9, 20).Value = 0.0001116986572225903

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Alcam"
$ws.Cells.Item(10, 3).Value = "L1cam"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.482215
$ws.Cells.Item(10, 8).Value = 1.446645
$ws.Cells.Item(10, 9).Value = 0.01120445187525872
$ws.Cells.Item(10, 10).Value = 0.01120445187525872
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.654706
$ws.Cells.Item(10, 14).Value = 22.964118
$ws.Cells.Item(10, 15).Value = 0.2505289128669849
$ws.Cells.Item(10, 16).Value = 0.2505289128669849
$ws.Cells.Item(10, 17).Value = 3.69121405379
$ws.Cells.Item(10, 18).Value = 33.22092648411
$ws.Cells.Item(10, 19).Value = 0.002807039147579018
$ws.Cells.Item(10, 20).Value = 0.002807039147579018

